$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-apply the AutoFilter on column A ("الاسم") so it only shows rows whose
# name is "حملة خرسان الثانية" (row 104). This both rewrites the stored
# <filters> list in the autoFilter definition and updates the hidden state
# of every data row to match the new criteria (rows that used to match the
# old filter values become hidden, and row 104 becomes visible).
$ws.Range("A1:O1033").AutoFilter(1, @("حملة خرسان الثانية"), 7)

# Move/record the active selection on the sheet to A104 (matches the new
# visible/filtered row), mirroring the <selection activeCell="A104" .../>
# recorded in the saved view state.
$ws.Range("A104").Select()
